$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 145 (Angeleno / 44266 entry),
# shifting the existing rows 145-156 down to 147-158.
$ws.Rows("145:146").Insert()

# Row 145: new "Black Amber" / "Especial" weekly entry
$ws.Range("A145").Value = 8
$ws.Range("B145").Value = "Terminal La Palmera de La Serena"
$ws.Range("C145").Value = "Coquimbo"
$ws.Range("D145").Value = 44578
$ws.Range("E145").Value = 4
$ws.Range("F145").Value = "Fruta"
$ws.Range("G145").Value = 100103
$ws.Range("H145").Value = "Frutos de hueso (carozo)"
$ws.Range("I145").Value = 100103002
$ws.Range("J145").Value = "Ciruela"
$ws.Range("K145").Value = "Black Amber"
$ws.Range("L145").Value = "Especial"
$ws.Range("M145").Value = 16
$ws.Range("N145").Value = 305000
$ws.Range("O145").Value = 310000
$ws.Range("P145").Value = 307500
$ws.Range("Q145").Value = "`$/bins (450 kilos)"
$ws.Range("R145").Value = "Región Metropolitana"
$ws.Range("S145").Value = 683
$ws.Range("T145").Value = 450

# Row 146: new "Black Amber" / "Primera" weekly entry
$ws.Range("A146").Value = 8
$ws.Range("B146").Value = "Terminal La Palmera de La Serena"
$ws.Range("C146").Value = "Coquimbo"
$ws.Range("D146").Value = 44578
$ws.Range("E146").Value = 4
$ws.Range("F146").Value = "Fruta"
$ws.Range("G146").Value = 100103
$ws.Range("H146").Value = "Frutos de hueso (carozo)"
$ws.Range("I146").Value = 100103002
$ws.Range("J146").Value = "Ciruela"
$ws.Range("K146").Value = "Black Amber"
$ws.Range("L146").Value = "Primera"
$ws.Range("M146").Value = 20
$ws.Range("N146").Value = 275000
$ws.Range("O146").Value = 280000
$ws.Range("P146").Value = 277500
$ws.Range("Q146").Value = "`$/bins (450 kilos)"
$ws.Range("R146").Value = "Región Metropolitana"
$ws.Range("S146").Value = 617
$ws.Range("T146").Value = 450

# Match the date-number format used by the rest of column D
$ws.Range("D145:D146").NumberFormat = $ws.Range("D147").NumberFormat()
